# Correct the "False Negatives" metric on the "Final Model - Metrics" slide.
#
# The slide's content placeholder holds four short paragraphs
# (Accuracy / Recall / False Negatives / AUC) as separate runs. We locate
# the paragraph that starts with "False Negatives:" and rewrite only that
# run's text, leaving every other run/paragraph (and their rPr formatting)
# untouched.

$p = $ppt.ActivePresentation

$oldText = "False Negatives: 216 / 7043 or 3%"
$newText = "False Negatives: 216 / 1761 or 12.3%"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)

    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)

        if (-not $shape.HasTextFrame) {
            continue
        }
        if (-not $shape.TextFrame.HasText) {
            continue
        }

        $textRange = $shape.TextFrame.TextRange
        $paragraphs = $textRange.Paragraphs()

        for ($paraIdx = 1; $paraIdx -le $paragraphs.Count; $paraIdx++) {
            $paragraph = $textRange.Paragraphs($paraIdx, 1)

            # Paragraphs(...).Text includes the trailing paragraph-mark
            # (carriage return) for every paragraph but the last one, so
            # trim it before comparing against the plain target string.
            $paraText = $paragraph.Text.TrimEnd([char]13)

            if ($paraText -eq $oldText) {
                $run = $paragraph.Runs(1, 1)
                $run.Text = $newText
            }
        }
    }
}
